$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("B2:G3").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "090019"
$newSheet.Range("C2").Value = "大成景恒混合A"
$newSheet.Range("D2").Value = "1.13"
$newSheet.Range("E2").Value = "93.98"
$newSheet.Range("F2").Value = "1.98"
$newSheet.Range("G2").Value = "0.0224"
$newSheet.Range("H2").Value = 4

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "006038"
$newSheet.Range("C3").Value = "大成景恒混合C"
$newSheet.Range("D3").Value = "0.45"
$newSheet.Range("E3").Value = "93.98"
$newSheet.Range("F3").Value = "1.98"
$newSheet.Range("G3").Value = "0.0089"
$newSheet.Range("H3").Value = 4

$newSheet.Range("B2:G3").ClearFormats()

$newSheet.Move($null, $summary)
